$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (report issue number & date range)
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value2  = "Volume 31   Number  18"
$ws.Cells.Item(9, 3).Value2  = "Report Covering the Week  4/29/2024  Through  5/5/2024"

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = -50
$ws.Range("N14").Value = -86.666666666666

# ---------------------------------------------------------------------------
# Row 15 - Rape (D15, E15 switch from numbers to text placeholders)
# ---------------------------------------------------------------------------
$ws.Range("C15").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -91.666666666666

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 84.615384615384
$ws.Range("I16").Value = 74
$ws.Range("J16").Value = 58
$ws.Range("K16").Value = 27.586206896551
$ws.Range("L16").Value = 34.545454545454
$ws.Range("M16").Value = -11.904761904761
$ws.Range("N16").Value = -72.794117647058

# ---------------------------------------------------------------------------
# Row 17 - Felonious Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 180
$ws.Range("F17").Value = 39
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 30
$ws.Range("I17").Value = 149
$ws.Range("J17").Value = 116
$ws.Range("K17").Value = 28.448275862069
$ws.Range("L17").Value = 53.608247422680
$ws.Range("M17").Value = 144.262295081967
$ws.Range("N17").Value = -21.164021164021

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 65
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 62.5
$ws.Range("L18").Value = 16.071428571428
$ws.Range("M18").Value = 261.111111111111
$ws.Range("N18").Value = -53.571428571428

# ---------------------------------------------------------------------------
# Row 19 - Grand Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -14.285714285714
$ws.Range("I19").Value = 124
$ws.Range("J19").Value = 145
$ws.Range("K19").Value = -14.482758620689
$ws.Range("L19").Value = -8.823529411764
$ws.Range("M19").Value = 77.142857142857
$ws.Range("N19").Value = -39.215686274509

# ---------------------------------------------------------------------------
# Row 20 - G.L.A. (C20, D20, E20 switch from text placeholders to numbers)
# ---------------------------------------------------------------------------
$ws.Range("F20").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 2
$ws.Range("F20").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 2
$ws.Range("H20").Copy($ws.Range("E20"))
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -40
$ws.Range("I20").Value = 20
$ws.Range("J20").Value = 36
$ws.Range("K20").Value = -44.444444444444
$ws.Range("L20").Value = -41.176470588235
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = -83.606557377049

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 54.545454545454
$ws.Range("F21").Value = 106
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = 19.101123595505
$ws.Range("I21").Value = 436
$ws.Range("J21").Value = 406
$ws.Range("K21").Value = 7.389162561576
$ws.Range("L21").Value = 12.082262210796
$ws.Range("M21").Value = 71.653543307086
$ws.Range("N21").Value = -54.865424430641

# ---------------------------------------------------------------------------
# Row 22 - Transit (F22 switches from number to text placeholder)
# ---------------------------------------------------------------------------
$ws.Range("C22").Copy($ws.Range("F22"))
$ws.Range("H22").Value = -100

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = 11
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 266.666666666667
$ws.Range("F23").Value = 32
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = 14.285714285714
$ws.Range("I23").Value = 126
$ws.Range("J23").Value = 121
$ws.Range("K23").Value = 4.132231404958
$ws.Range("L23").Value = 4.132231404958
$ws.Range("M23").Value = 88.059701492537

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -23.809523809523
$ws.Range("F24").Value = 68
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = -20.930232558139
$ws.Range("I24").Value = 274
$ws.Range("J24").Value = 304
$ws.Range("K24").Value = -9.868421052631
$ws.Range("L24").Value = 12.295081967213
$ws.Range("M24").Value = 35.643564356435

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -55.555555555555
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = -50
$ws.Range("I25").Value = 52
$ws.Range("J25").Value = 75
$ws.Range("K25").Value = -30.666666666666
$ws.Range("L25").Value = 1.960784313725

# ---------------------------------------------------------------------------
# Row 26 - Misdemeanor Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 250
$ws.Range("F26").Value = 66
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = 60.975609756097
$ws.Range("I26").Value = 242
$ws.Range("J26").Value = 186
$ws.Range("K26").Value = 30.107526881720
$ws.Range("L26").Value = 29.411764705882
$ws.Range("M26").Value = 7.555555555555

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("G27").Value = 2
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = -45.454545454545

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 20
$ws.Range("K28").Value = 17.647058823529
$ws.Range("L28").Value = 11.111111111111

# ---------------------------------------------------------------------------
# Row 29 - Shooting Victims
# ---------------------------------------------------------------------------
$ws.Range("M29").Value = -78.571428571428
$ws.Range("N29").Value = -91.428571428571

# ---------------------------------------------------------------------------
# Row 30 - Shooting Incidents
# ---------------------------------------------------------------------------
$ws.Range("M30").Value = -83.333333333333
$ws.Range("N30").Value = -93.75
